# "added login features scenarios"
# Replace the code samples shown in cells A7 and A9 of the "PythonCode" sheet
# with newer / refactored versions of the same algorithms. Excel will append
# the two new strings to the shared-string table (becoming entries 16 and 17)
# while cells A6/A8 keep pointing at the original (now otherwise unused)
# shared strings 9 and 11.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PythonCode")

$newMaxConsecutiveOnes = @"
def findMaxConsecutiveOnes(nums) :
max_count = 0
current_count = 0
for num in nums:
if num == 1:
current_count += 1
max_count = max(max_count, current_count)
\b
\b
else:
current_count = 0
#\b\b\b
\b
\b
\b
\b
return max_count
"@

$newFindNumbers = 'def findNumbers(nums):return sum(len(str(num)) % 2 == 0 for num in nums)'

$ws.Range("A7").Value = $newMaxConsecutiveOnes
$ws.Rows.Item(7).AutoFit()
$ws.Range("A9").Value = $newFindNumbers

[void]$ws.Range("A9").Select()
